# Update public EPEX Spot prices workbook with the latest daily data point.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": append a new date column (06-sep) after the last one (CF)
# ---------------------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Copy the header cell's formatting (bold, bordered, centered) onto the new
# header cell before writing its text, then set the label.
$wsSpot.Range("CF1").Copy()
$wsSpot.Range("CG1").PasteSpecial(-4122)
$wsSpot.Range("CG1").Value = "06-sep"

# Hourly prices for 06-sep.
$wsSpot.Range("CG2").Value = 74.64
$wsSpot.Range("CG3").Value = 57.8
$wsSpot.Range("CG4").Value = 41.35
$wsSpot.Range("CG5").Value = 20.66
$wsSpot.Range("CG6").Value = 20
$wsSpot.Range("CG7").Value = 25.38
$wsSpot.Range("CG8").Value = 25.49
$wsSpot.Range("CG9").Value = 25.86
$wsSpot.Range("CG10").Value = 29.28
$wsSpot.Range("CG11").Value = 15.23
$wsSpot.Range("CG12").Value = 1.72
$wsSpot.Range("CG13").Value = 0
$wsSpot.Range("CG14").Value = -0.02
$wsSpot.Range("CG15").Value = -0.86
$wsSpot.Range("CG16").Value = -0.99
$wsSpot.Range("CG17").Value = -0.01
$wsSpot.Range("CG18").Value = -0.01
$wsSpot.Range("CG19").Value = 0
$wsSpot.Range("CG20").Value = 17.42
$wsSpot.Range("CG21").Value = 24.37
$wsSpot.Range("CG22").Value = 42
$wsSpot.Range("CG23").Value = 19.01
$wsSpot.Range("CG24").Value = 36.26
$wsSpot.Range("CG25").Value = 27.99

# ---------------------------------------------------------------------------
# Sheet "Gaz": append the 2025-09-04 daily price row
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

# Force the date column to text so "2025-09-04" is stored as a literal
# string (matching every other row) instead of being auto-converted to a
# date serial, then restore the plain/default cell formatting used by the
# rest of the column.
$wsGaz.Range("A82").NumberFormat = "@"
$wsGaz.Range("A82").Value = "2025-09-04"
$wsGaz.Range("A81").Copy()
$wsGaz.Range("A82").PasteSpecial(-4122)

$wsGaz.Range("B82").Value = 31.5

# ---------------------------------------------------------------------------
# Sheet "CO2": append the 2025-09-04 daily price row
# ---------------------------------------------------------------------------
$wsCO2 = $wb.Worksheets.Item("CO2")

$wsCO2.Range("A82").NumberFormat = "@"
$wsCO2.Range("A82").Value = "2025-09-04"
$wsCO2.Range("A81").Copy()
$wsCO2.Range("A82").PasteSpecial(-4122)

$wsCO2.Range("B82").Value = 74.90000000000001
